$wb = $excel.ActiveWorkbook

# Step 1: move the selection on sheet 'PAR RAPPORT A UN POINT' to F2 BEFORE adding
# any new sheets, so it does not end up being the tabSelected sheet at the end.
$ws4 = $wb.Worksheets.Item("PAR RAPPORT À UN POINT")
$ws4.Range("F2").Select()

# Step 2: append the two new sheets, in order, after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws5 = $wb.Worksheets.Add($null, $lastSheet)
$ws5.Name = "DISTANCE ET PROXIMITÉ"
$ws6 = $wb.Worksheets.Add($null, $ws5)
$ws6.Name = "DIRECTION ET MOUVEMENT"

# Step 3: copy the header row formatting (bold/fill/border) from an existing sheet
$headerSrc = $wb.Worksheets.Item("GÉNÉRALITÉS")
$headerSrc.Range("A1:D1").Copy()
$ws5.Range("A1:D1").PasteSpecial(-4122)
$headerSrc.Range("A1:D1").Copy()
$ws6.Range("A1:D1").PasteSpecial(-4122)

# Populate 'DISTANCE ET PROXIMITÉ'
$ws5.Cells.Item(1,1).Value = "Mot français"
$ws5.Cells.Item(1,2).Value = "Grammaire"
$ws5.Cells.Item(1,3).Value = "Prononciation"
$ws5.Cells.Item(1,4).Value = "Signification en tchèque"
$ws5.Cells.Item(2,1).Value = 'aux alentours'
$ws5.Cells.Item(2,2).Value = 'adv'
$ws5.Cells.Item(2,3).Value = 'o.zala~tu:r'
$ws5.Cells.Item(2,4).Value = 'v okolí'
$ws5.Cells.Item(3,1).Value = 'contre'
$ws5.Cells.Item(3,2).Value = 'prép'
$ws5.Cells.Item(3,3).Value = 'ko~:tr'
$ws5.Cells.Item(3,4).Value = 'proti, naproti'
$ws5.Cells.Item(4,1).Value = 'tout contre'
$ws5.Cells.Item(4,2).Value = 'adv'
$ws5.Cells.Item(4,3).Value = 'tu ko~:tr'
$ws5.Cells.Item(4,4).Value = 'zcela blízko'
$ws5.Cells.Item(5,1).Value = 'à côté (de)'
$ws5.Cells.Item(5,3).Value = 'a ko.te. d@'
$ws5.Cells.Item(5,4).Value = 'vedle (čeho), mimo'
$ws5.Cells.Item(6,1).Value = 'au côté de'
$ws5.Cells.Item(6,2).Value = 'prép'
$ws5.Cells.Item(6,3).Value = 'o. ko.te. d@'
$ws5.Cells.Item(6,4).Value = 'po boku čeho'
$ws5.Cells.Item(7,1).Value = 'du côté (de)'
$ws5.Cells.Item(7,3).Value = 'dü ko.te. d@'
$ws5.Cells.Item(7,4).Value = 'vedle (čeho)'
$ws5.Cells.Item(8,1).Value = 'à l''écart (de)'
$ws5.Cells.Item(8,3).Value = 'a le.ka:r d@'
$ws5.Cells.Item(8,4).Value = 'stranou (čeho)'
$ws5.Cells.Item(9,1).Value = 'dans les environs (de)'
$ws5.Cells.Item(9,3).Value = 'da~ le.za~viro~ d@'
$ws5.Cells.Item(9,4).Value = 'v okolí (čeho)'
$ws5.Cells.Item(10,1).Value = 'loin (de)'
$ws5.Cells.Item(10,3).Value = 'lu^e~ d@'
$ws5.Cells.Item(10,4).Value = 'daleko (od čeho)'
$ws5.Cells.Item(11,1).Value = 'au loin'
$ws5.Cells.Item(11,2).Value = 'adv'
$ws5.Cells.Item(11,3).Value = 'o. lu^e~'
$ws5.Cells.Item(11,4).Value = 'v dálce, daleko'
$ws5.Cells.Item(12,1).Value = 'un peu plus loin'
$ws5.Cells.Item(12,2).Value = 'adv'
$ws5.Cells.Item(12,3).Value = 'ö~ pö plü lu^e~'
$ws5.Cells.Item(12,4).Value = 'trochu dál'
$ws5.Cells.Item(13,1).Value = 'à mi-chemin'
$ws5.Cells.Item(13,2).Value = 'adv'
$ws5.Cells.Item(13,3).Value = 'a miš@me~'
$ws5.Cells.Item(13,4).Value = 'v půli cesty, na půl cesty'
$ws5.Cells.Item(14,1).Value = 'à peu de distance de'
$ws5.Cells.Item(14,2).Value = 'prép'
$ws5.Cells.Item(14,3).Value = 'a pö d@ dista~:s d@'
$ws5.Cells.Item(14,4).Value = 'na vzdálenost o málo kratší než'
$ws5.Cells.Item(15,1).Value = 'près (de)'
$ws5.Cells.Item(15,3).Value = 'pre d@'
$ws5.Cells.Item(15,4).Value = 'blízko (čeho)'
$ws5.Cells.Item(16,1).Value = 'tout près'
$ws5.Cells.Item(16,2).Value = 'adv'
$ws5.Cells.Item(16,3).Value = 'tu pre'
$ws5.Cells.Item(16,4).Value = 'zcela blízko'
$ws5.Cells.Item(17,1).Value = 'à proximité de'
$ws5.Cells.Item(17,2).Value = 'prép'
$ws5.Cells.Item(17,3).Value = 'a proksimite. d@'
$ws5.Cells.Item(17,4).Value = 'v blízkosti čeho'

# Populate 'DIRECTION ET MOUVEMENT'
$ws6.Cells.Item(1,1).Value = "Mot français"
$ws6.Cells.Item(1,2).Value = "Grammaire"
$ws6.Cells.Item(1,3).Value = "Prononciation"
$ws6.Cells.Item(1,4).Value = "Signification en tchèque"
$ws6.Cells.Item(2,1).Value = 'd''un bout à l''autre'
$ws6.Cells.Item(2,3).Value = 'dö~ bu a lo.tr'
$ws6.Cells.Item(2,4).Value = 'z jednoho konce na druhý'
$ws6.Cells.Item(3,1).Value = 'en chemin'
$ws6.Cells.Item(3,2).Value = 'adv'
$ws6.Cells.Item(3,3).Value = 'a~ š@me~'
$ws6.Cells.Item(3,4).Value = 'na cestě, po cestě, cestou, mezitím'
$ws6.Cells.Item(4,1).Value = 'en direction de'
$ws6.Cells.Item(4,2).Value = 'prép'
$ws6.Cells.Item(4,3).Value = 'a~ direksjo~ d@'
$ws6.Cells.Item(4,4).Value = 'směrem na co, směrem k čemu, ve směru čeho'
$ws6.Cells.Item(5,1).Value = 'dans la direction de'
$ws6.Cells.Item(5,2).Value = 'prép'
$ws6.Cells.Item(5,3).Value = 'da~ la direksjo~ d@'
$ws6.Cells.Item(5,4).Value = 'směrem na co'
$ws6.Cells.Item(6,1).Value = 'tout droit'
$ws6.Cells.Item(6,2).Value = 'adv'
$ws6.Cells.Item(6,3).Value = 'tu dru^a'
$ws6.Cells.Item(6,4).Value = 'přímo, pořád rovně'
$ws6.Cells.Item(7,1).Value = 'jusqu''à'
$ws6.Cells.Item(7,2).Value = 'prép'
$ws6.Cells.Item(7,3).Value = 'žüska'
$ws6.Cells.Item(7,4).Value = 'až k (čemu)'
$ws6.Cells.Item(8,1).Value = 'par'
$ws6.Cells.Item(8,2).Value = 'prép'
$ws6.Cells.Item(8,3).Value = 'par'
$ws6.Cells.Item(8,4).Value = 'přes'
$ws6.Cells.Item(9,1).Value = 'par là'
$ws6.Cells.Item(9,2).Value = 'adv'
$ws6.Cells.Item(9,3).Value = 'par la'
$ws6.Cells.Item(9,4).Value = 'tudy'
$ws6.Cells.Item(10,1).Value = 'à travers'
$ws6.Cells.Item(10,2).Value = 'prép'
$ws6.Cells.Item(10,3).Value = 'a trave:r'
$ws6.Cells.Item(10,4).Value = 'napříč, přes co'
$ws6.Cells.Item(11,1).Value = 'sur la route de'
$ws6.Cells.Item(11,3).Value = 'sür la rut d@'
$ws6.Cells.Item(11,4).Value = 'na cestě k čemu'
$ws6.Cells.Item(12,1).Value = 'vers'
$ws6.Cells.Item(12,2).Value = 'prép'
$ws6.Cells.Item(12,3).Value = 've:r'
$ws6.Cells.Item(12,4).Value = 'směrem k (čemu)'

# Step 4: set the F2 JSON-export helper formula on each new sheet (same pattern used
# throughout the rest of the workbook).
$ws5.Cells.Item(2,6).Formula = '= "{ ""foreign"": """ & A2 & """, ""grammar"": """ & B2 & """, ""pronunciation"": """ & C2 & """, ""meaning"": """ & D2 & """ },"'
$ws6.Cells.Item(2,6).Formula = '= "{ ""foreign"": """ & A2 & """, ""grammar"": """ & B2 & """, ""pronunciation"": """ & C2 & """, ""meaning"": """ & D2 & """ },"'

# Step 5: final selections to match the target sheet views.
$ws5.Range("F2").Select()
$ws6.Range("A13").Select()
